$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): columns C, D, E get new labels
$ws.Range("C1").Value2 = "prediction"
$ws.Range("D1").Value2 = "rejection-f"
$ws.Range("E1").Value2 = "max"

# Determine last used row in column A (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    # Column C becomes the same text value as column D (genus prediction)
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal

    # Column E becomes numeric 1
    $ws.Cells.Item($r, 5).Value2 = 1
}
